$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($sheet, [string]$addr, [string]$text)
    # Prefix with an apostrophe so Excel stores numeric-looking strings
    # as text (matching the source file's inline-string / text cells)
    # instead of silently converting them to numbers.
    $sheet.Range($addr).Value = "'" + $text
}

# Price (column D) updates
Set-TextValue $ws "D2"  "283.60"
Set-TextValue $ws "D3"  "20.84"
Set-TextValue $ws "D4"  "6.197"
Set-TextValue $ws "D5"  "0.06190"
Set-TextValue $ws "D6"  "3.582"
Set-TextValue $ws "D8"  "1.485"
Set-TextValue $ws "D9"  "0.8172"
Set-TextValue $ws "D10" "0.01390"
Set-TextValue $ws "D11" "0.1649"
Set-TextValue $ws "D12" "0.08372"
Set-TextValue $ws "D13" "0.03670"
Set-TextValue $ws "D14" "0.03133"
Set-TextValue $ws "D15" "0.09133"
Set-TextValue $ws "D16" "3.726"
Set-TextValue $ws "D17" "0.001639"
Set-TextValue $ws "D18" "0.04667"
Set-TextValue $ws "D19" "0.006466"
Set-TextValue $ws "D20" "0.006193"
Set-TextValue $ws "D22" "0.0001500"
Set-TextValue $ws "D23" "3.799"
Set-TextValue $ws "D26" "0.1220"

# Rows 41-43 got re-ranked: coins shifted up one slot (row43 -> row41 slot
# content-wise cycles), each with its own fresh price/volume figure.
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.007067"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1106"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.003509"
$ws.Range("E43").Value = "42CEJICEJI"

# Remaining price (column D) updates
Set-TextValue $ws "D44" "0.01147"
Set-TextValue $ws "D45" "0.00006440"
Set-TextValue $ws "D47" "0.8402"
Set-TextValue $ws "D50" "0.01240"
